$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Date cells: update to 2024-04-29 (Excel serial 45411) ---
$ws.Range("B6").Value = 45411
$ws.Range("B12").Value = 45411
$ws.Range("F12").Value = 45411

# --- Numeric cells: replace amounts with 1 ---
$ws.Range("F20").Value = 1
$ws.Range("F32").Value = 1

# --- Text cells: replace with literal text "1" (not a number) ---
# Use a scratch cell formatted as Text so the value "1" is typed in as a
# genuine string, then copy/paste-special (values only) into every target
# cell so each destination keeps its own existing style/format untouched.
$scratch = $ws.Range("Z100")
$scratch.NumberFormat = "@"
$scratch.Value = "1"
$scratch.Copy()

$textCells = @("B7", "B8", "F8", "B9", "B11", "B13", "A20", "B20", "C20", "D20", "E20", "A40", "E40", "B45", "B46")
foreach ($addr in $textCells) {
    $ws.Range($addr).PasteSpecial(-4163)
}

$scratch.Clear()
$excel.CutCopyMode = $false
